# Edit: (1) change the table style on slide 5's table to the new style
#       (2) repaint the deck's (single) live theme colour scheme with the
#           "Office Theme" palette, matching the swapped ppt/theme/theme1.xml
#           colour values from the target OOXML.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{5A8EFEE3-AFD7-4136-8369-C00D17920B8C}")
    }
}

# --- 2. Theme colours ------------------------------------------------------
# msoThemeColor order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeRGB = @{
    1  = 0           # dk1      000000
    2  = 16777215    # lt1      FFFFFF
    3  = 6968388      # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845    # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308    # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797    # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeRGB[$i]
}

$p.Save()
